$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")

# Update column A values: "Us Batw 8 47".."Us Batw 8 51" -> "Us Batw 8 58".."Us Batw 8 62"
$ws.Cells.Item(1, 1).Value = "Us Batw 8 58"
$ws.Cells.Item(2, 1).Value = "Us Batw 8 59"
$ws.Cells.Item(3, 1).Value = "Us Batw 8 60"
$ws.Cells.Item(4, 1).Value = "Us Batw 8 61"
$ws.Cells.Item(5, 1).Value = "Us Batw 8 62"

# Update column B values: "ubw47".."ubw51" -> "ubw58".."ubw62"
$ws.Cells.Item(1, 2).Value = "ubw58"
$ws.Cells.Item(2, 2).Value = "ubw59"
$ws.Cells.Item(3, 2).Value = "ubw60"
$ws.Cells.Item(4, 2).Value = "ubw61"
$ws.Cells.Item(5, 2).Value = "ubw62"

# Update the selection on this sheet to A1:B5 (no active cell at C7)
$ws.Range("A1:B5").Select()
